$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.789.11'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -4.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.146.16'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -8.70%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.58'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.59'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.71%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -3.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.144.87'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -8.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.123'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -7.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.52'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -6.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.391'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -6.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.699.35'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -8.46%  '
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.90'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -9.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.801.75'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000161'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -7.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.153.09'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -7.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.69'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.83'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -7.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '350.09'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.14'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -6.59%  '
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.88'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -7.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.496'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -7.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000114'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -11.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.44'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.92%  '
$ws.Range("E28").Value = '  -1.80%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.87'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -6.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.39'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -8.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.79'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -7.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.19'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -6.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.55'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -7.48%  '
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '153.39'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -6.14%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.41'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -9.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.815'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -7.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.17'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -5.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.67'
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.628.44'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.41'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -7.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.12'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -8.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.12'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.95'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -8.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0644'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -7.46%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.45'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -7.83%  '
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '317.44'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -6.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0268'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -6.75%  '
$ws.Range("E50").Value = '  -4.58%  '
$ws.Range("E51").Value = '  +0.06%  '
